$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-12-25 Monday" "2023-12-26 Tuesday"
Replace-Text "95×84=" "55×49="
Replace-Text "97×92=" "29×21="
Replace-Text "98×58=" "46×16="
Replace-Text "61×79=" "95×15="
Replace-Text "30×85=" "86×37="
Replace-Text "25×29=" "57×38="
Replace-Text "20×99=" "26×68="
Replace-Text "67×45=" "98×64="
Replace-Text "32×28=" "84×94="
Replace-Text "24×84=" "33×76="
Replace-Text "39×90=" "53×35="
Replace-Text "41×36=" "35×14="
Replace-Text "86×86=" "30×71="
Replace-Text "74×85=" "74×60="
Replace-Text "86×68=" "73×96="
Replace-Text "63×47=" "40×69="
Replace-Text "13×68=" "19×90="
Replace-Text "91×46=" "92×72="
Replace-Text "89×50=" "19×45="
Replace-Text "61×76=" "92×45="
Replace-Text "74×73=" "34×13="
Replace-Text "57×48=" "89×35="
Replace-Text "70×22=" "25×78="
Replace-Text "73×34=" "42×15="
Replace-Text "27×29=" "50×57="
